# Append new scrape results: 2025-10-21 12:39 JST
# Target sheet: "ランサーズ" (the 1st worksheet / listing sheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTs = "2025-10-21 12:38:58"

# ------------------------------------------------------------------
# 1. Widen column B (タイトル / title) from 45 to 50 characters.
#    The engine stores width in "raw" OOXML character units that are
#    derived from the ColumnWidth (Calibri 11 char units) through a
#    pixel-rounding step, so we feed it a value that lands exactly on
#    raw width 50 once rounded.
# ------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 49.17

# ------------------------------------------------------------------
# 2. A brand-new listing was scraped and inserted right after the
#    existing top 4 rows (rows 2-5 keep their data, only row 6 onward
#    shifts down by one row).
# ------------------------------------------------------------------
$ws.Rows.Item(6).Insert()

# ------------------------------------------------------------------
# 3. Refresh the "取得日時" (fetched-at) timestamp for every row that
#    belongs to this scrape run (the previously existing rows 2-5 and
#    the ones that shifted down to 7-9, plus the new rows 6, 10, 11).
# ------------------------------------------------------------------
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTs
}

# ------------------------------------------------------------------
# 4. Populate the new row 6 with the newly scraped listing.
# ------------------------------------------------------------------
$ws.Range("B6").Value = "【急募】MT4/MT5用FX自動売買システムの開発者募集"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5417377"
$ws.Range("G6").Value = 83
$ws.Range("H6").Value = "◆開発"

# ------------------------------------------------------------------
# 5. Append two more freshly scraped listings at the end of the table
#    (rows 10 and 11). The last one has no skill-summary tag.
# ------------------------------------------------------------------
$ws.Range("A10").Value = $newTs
$ws.Range("B10").Value = "【Webarena suiteX/DNS】ドメイン設定変更によるウェブサイト分割とサイト切り替え"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5417544"
$ws.Range("G10").Value = 30
$ws.Range("H10").Value = "◇サイト"

$ws.Range("A11").Value = $newTs
$ws.Range("B11").Value = "【データ加工のプロ募集】施設情報データの修正・整備依頼"
$ws.Range("C11").Value = "システム開発"
$ws.Range("D11").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E11").Value = "期限情報なし"
$ws.Range("F11").Value = "https://www.lancers.jp/work/detail/5417622"
$ws.Range("G11").Value = 10

# ------------------------------------------------------------------
# 6. Fix up the URL hyperlinks. Inserting the row above only shifted
#    the cell *data* down, not the hyperlink anchors, so the safest
#    way to end up with correct, consistent hyperlinks for every row
#    is to clear all of them and recreate them against the final
#    (post-shift) layout, then restore the "Hyperlink" cell style.
# ------------------------------------------------------------------
$ws.Cells.Hyperlinks.Delete()

$urls = @{
    2 = "https://www.lancers.jp/work/detail/5415908"
    3 = "https://www.lancers.jp/work/detail/5417295"
    4 = "https://www.lancers.jp/work/detail/5408664"
    5 = "https://www.lancers.jp/work/detail/5417433"
    6 = "https://www.lancers.jp/work/detail/5417377"
    7 = "https://www.lancers.jp/work/detail/5371747"
    8 = "https://www.lancers.jp/work/detail/5417267"
    9 = "https://www.lancers.jp/work/detail/5417308"
    10 = "https://www.lancers.jp/work/detail/5417544"
    11 = "https://www.lancers.jp/work/detail/5417622"
}

for ($r = 2; $r -le 11; $r++) {
    $cell = $ws.Range("F" + $r)
    $ws.Hyperlinks.Add($cell, $urls[$r])
    $cell.Style = "Hyperlink"
}
